$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4504.627264712643
$ws.Range("B3").Value = 4633.33839461559
$ws.Range("B4").Value = 4759.722242835674
$ws.Range("B5").Value = 4883.392489427923
$ws.Range("B6").Value = 5003.962922665383
$ws.Range("B7").Value = 5121.049412950785
$ws.Range("B8").Value = 5234.271923436626
$ws.Range("B9").Value = 5343.256542883342
$ws.Range("B10").Value = 5447.637525739176
$ws.Range("B11").Value = 5547.059324019701
$ws.Range("B12").Value = 5641.178595307084
$ws.Range("B13").Value = 5729.666171085652
$ws.Range("B14").Value = 5812.208969685085
$ws.Range("B15").Value = 5888.511838318149
$ws.Range("B16").Value = 5958.299309076966
$ws.Range("B17").Value = 6021.317254287688
$ws.Range("B18").Value = 6077.334427315947
$ws.Range("B19").Value = 6126.143875757166
$ws.Range("B20").Value = 6167.564214930656
$ws.Range("B21").Value = 6201.440750713115
$ws.Range("B22").Value = 6227.64644198572
$ws.Range("B23").Value = 6246.082694315713
$ws.Range("B24").Value = 6256.679977933589
$ws.Range("B25").Value = 6259.39826458554
$ws.Range("B26").Value = 6254.227279420029
$ws.Range("B27").Value = 6241.186565690459
$ws.Range("B28").Value = 6220.325361703711
$ws.Range("B29").Value = 6191.722291098811
$ws.Range("B30").Value = 6155.484869182075
$ws.Range("B31").Value = 6111.748829656309
$ws.Range("B32").Value = 6060.677277644281
$ws.Range("B33").Value = 6002.459676402021
$ws.Range("B34").Value = 5937.310676530894
$ws.Range("B35").Value = 5865.468797811419
$ws.Range("B36").Value = 5787.194974984221
$ws.Range("B37").Value = 5702.770979880082
$ws.Range("B38").Value = 5612.497733242381
$ws.Range("B39").Value = 5516.693520381373
$ws.Range("B40").Value = 5415.692125443871
$ws.Range("B41").Value = 5309.840899568959
$ws.Range("B42").Value = 5199.498778527
$ws.Range("B43").Value = 5085.034265604899
$ws.Range("B44").Value = 4966.823395505629
$ws.Range("B45").Value = 4845.247694878099
$ws.Range("B46").Value = 4720.692154788529
$ws.Range("B47").Value = 4593.543229994942
$ws.Range("B48").Value = 4464.186879297748
$ws.Range("B49").Value = 4333.006660524843
$ws.Range("B50").Value = 4200.381892877557
$ws.Range("B51").Value = 4066.685898428344
$ws.Range("B52").Value = 3932.284333534301
$ws.Range("B53").Value = 3797.533619828217
$ws.Range("B54").Value = 3662.77948328296
$ws.Range("B55").Value = 3528.355608633436
$ws.Range("B56").Value = 3394.582415195665
$ws.Range("B57").Value = 3261.765958861217
$ws.Range("B58").Value = 3130.196963781129
$ws.Range("B59").Value = 3000.149986001161
$ws.Range("B60").Value = 2871.882710082897
$ws.Range("B61").Value = 2745.635378555737
$ws.Range("B62").Value = 2621.630352904689
$ws.Range("B63").Value = 2500.071803719092
$ws.Range("B64").Value = 2381.14552661683
$ws.Range("B65").Value = 2265.018879625954
$ws.Range("B66").Value = 2151.840836857543
$ws.Range("B67").Value = 2041.742152545537
$ws.Range("B68").Value = 1934.835628865538
$ws.Range("B69").Value = 1831.216480377853
$ws.Range("B70").Value = 1730.962787471794
$ws.Range("B71").Value = 1634.136030818801
$ws.Range("B72").Value = 1540.781698570106
$ws.Range("B73").Value = 1450.929957858329
$ws.Range("B74").Value = 1364.596382078415

$ws.Range("A75:B77").ClearContents()
